$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.86194475070221
$ws.Range("C2").Value = 8.765511308359851
$ws.Range("D2").Value = 3.662496967612788
$ws.Range("F2").Value = 21.9602068832595
$ws.Range("G2").Value = 25.67250057052817
$ws.Range("H2").Value = 12.8474512730467
$ws.Range("I2").Value = 17.96781131422806
$ws.Range("M2").Value = 19.96712029853189
$ws.Range("N2").Value = 16.4315220332652
$ws.Range("B3").Value = 10.31135853239955
$ws.Range("C3").Value = 8.303997148163383
$ws.Range("D3").Value = 3.647148156859569
$ws.Range("F3").Value = 21.77557571008817
$ws.Range("G3").Value = 25.24395389201186
$ws.Range("H3").Value = 12.85720858881789
$ws.Range("I3").Value = 18.00954440375927
$ws.Range("M3").Value = 19.35943292105068
$ws.Range("N3").Value = 16.50351005695181
$ws.Range("B4").Value = 9.959174187497171
$ws.Range("C4").Value = 8.005361946026103
$ws.Range("D4").Value = 3.637579773793059
$ws.Range("F4").Value = 21.66951450099396
$ws.Range("G4").Value = 24.98880842876541
$ws.Range("H4").Value = 12.86640559219111
$ws.Range("I4").Value = 18.04033248779136
$ws.Range("M4").Value = 18.9846351146308
$ws.Range("N4").Value = 16.54956241880686
$ws.Range("B5").Value = 9.812282175279332
$ws.Range("C5").Value = 7.879887664047845
$ws.Range("D5").Value = 3.633645650873091
$ws.Range("F5").Value = 21.62817464366716
$ws.Range("G5").Value = 24.8870017372267
$ws.Range("H5").Value = 12.87095776561364
$ws.Range("I5").Value = 18.0541707453555
$ws.Range("M5").Value = 18.83173848483743
$ws.Range("N5").Value = 16.56879678073527
$ws.Range("B6").Value = 9.78769285600568
$ws.Range("C6").Value = 7.858826586289217
$ws.Range("D6").Value = 3.632990321707335
$ws.Range("F6").Value = 21.62142501298327
$ws.Range("G6").Value = 24.87023226298064
$ws.Range("H6").Value = 12.87176216798174
$ws.Range("I6").Value = 18.05654640148146
$ws.Range("M6").Value = 18.80634802907286
$ws.Range("N6").Value = 16.57201893808112
$ws.Range("B7").Value = 9.95720655565945
$ws.Range("C7").Value = 8.003684964095275
$ws.Range("D7").Value = 3.637526856541805
$ws.Range("F7").Value = 21.6689493057907
$ws.Range("G7").Value = 24.98742644655866
$ws.Range("H7").Value = 12.86646373051163
$ws.Range("I7").Value = 18.04051389380002
$ws.Range("M7").Value = 18.98257340903521
$ws.Range("N7").Value = 16.54981992394851
$ws.Range("B8").Value = 10.67513094268495
$ws.Range("C8").Value = 8.609581445087548
$ws.Range("D8").Value = 3.657235404608029
$ws.Range("F8").Value = 21.89505831422094
$ws.Range("G8").Value = 25.5231777987405
$ws.Range("H8").Value = 12.85014919573495
$ws.Range("I8").Value = 17.98112527935036
$ws.Range("M8").Value = 19.75808783630529
$ws.Range("N8").Value = 16.45596051034463
$ws.Range("B9").Value = 11.9645468455595
$ws.Range("C9").Value = 9.67458477411569
$ws.Range("D9").Value = 3.694683063391323
$ws.Range("F9").Value = 22.39424535836159
$ws.Range("G9").Value = 26.62991791227055
$ws.Range("H9").Value = 12.84366143683344
$ws.Range("I9").Value = 17.9059148754663
$ws.Range("M9").Value = 21.25498684841351
$ws.Range("N9").Value = 16.28649532637317
$ws.Range("B10").Value = 12.8328055026891
$ws.Range("C10").Value = 10.3799040943613
$ws.Range("D10").Value = 3.721388896259307
$ws.Range("F10").Value = 22.79186101184996
$ws.Range("G10").Value = 27.46746456653022
$ws.Range("H10").Value = 12.85451200200682
$ws.Range("I10").Value = 17.87615839372357
$ws.Range("M10").Value = 22.32717912165167
$ws.Range("N10").Value = 16.17074648944493
$ws.Range("B11").Value = 13.20957282453545
$ws.Range("C11").Value = 10.68374267348345
$ws.Range("D11").Value = 3.733347058542964
$ws.Range("F11").Value = 22.97871685065532
$ws.Range("G11").Value = 27.85172492080215
$ws.Range("H11").Value = 12.86284535934976
$ws.Range("I11").Value = 17.86822425365484
$ws.Range("M11").Value = 22.80655757884262
$ws.Range("N11").Value = 16.11996080741972
$ws.Range("B12").Value = 13.35964950910524
$ws.Range("C12").Value = 10.79633565108116
$ws.Range("D12").Value = 3.737846568500422
$ws.Range("F12").Value = 23.05026950869142
$ws.Range("G12").Value = 27.99753482191915
$ws.Range("H12").Value = 12.8664892337072
$ws.Range("I12").Value = 17.86602990145908
$ws.Range("M12").Value = 22.9867068596973
$ws.Range("N12").Value = 16.10099605299884
$ws.Range("B13").Value = 13.32601583011346
$ws.Range("C13").Value = 10.77219655720821
$ws.Range("D13").Value = 3.736878820986149
$ws.Range("F13").Value = 23.0348250869696
$ws.Range("G13").Value = 27.96612140089112
$ws.Range("H13").Value = 12.86568275676987
$ws.Range("I13").Value = 17.86646639719687
$ws.Range("M13").Value = 22.9479725468653
$ws.Range("N13").Value = 16.10506862472901
$ws.Range("B14").Value = 13.22114379574811
$ws.Range("C14").Value = 10.69305523254827
$ws.Range("D14").Value = 3.733717820843363
$ws.Range("F14").Value = 22.98458798896237
$ws.Range("G14").Value = 27.86371575156362
$ws.Range("H14").Value = 12.86313536205534
$ws.Range("I14").Value = 17.86802746267695
$ws.Range("M14").Value = 22.82140719927078
$ws.Range("N14").Value = 16.1183952329475
$ws.Range("B15").Value = 13.16052693247204
$ws.Range("C15").Value = 10.64425755267172
$ws.Range("D15").Value = 3.731777826751356
$ws.Range("F15").Value = 22.95391780806615
$ws.Range("G15").Value = 27.80102331359253
$ws.Range("H15").Value = 12.86163857154002
$ws.Range("I15").Value = 17.86908928892386
$ws.Range("M15").Value = 22.74369746045302
$ws.Range("N15").Value = 16.12659283301489
$ws.Range("B16").Value = 12.80781028857369
$ws.Range("C16").Value = 10.35970355007973
$ws.Range("D16").Value = 3.72060344995064
$ws.Range("F16").Value = 22.77976429996504
$ws.Range("G16").Value = 27.44240317749518
$ws.Range("H16").Value = 12.85403574390549
$ws.Range("I16").Value = 17.87679004055222
$ws.Range("M16").Value = 22.29566663936975
$ws.Range("N16").Value = 16.17410288506888
$ws.Range("B17").Value = 12.58671325154791
$ws.Range("C17").Value = 10.18076704996328
$ws.Range("D17").Value = 3.713698553020992
$ws.Range("F17").Value = 22.67441205576068
$ws.Range("G17").Value = 27.22311326685475
$ws.Range("H17").Value = 12.85024171715711
$ws.Range("I17").Value = 17.88295247544594
$ws.Range("M17").Value = 22.01853903406867
$ws.Range("N17").Value = 16.20372598878965
$ws.Range("B18").Value = 12.45783346470129
$ws.Range("C18").Value = 10.07624669945778
$ws.Range("D18").Value = 3.709709195549587
$ws.Range("F18").Value = 22.61438297993763
$ws.Range("G18").Value = 27.09729986132039
$ws.Range("H18").Value = 12.84837931912256
$ws.Range("I18").Value = 17.88702401710281
$ws.Range("M18").Value = 21.85836448079736
$ws.Range("N18").Value = 16.22094046299822
$ws.Range("B19").Value = 12.41390550209761
$ws.Range("C19").Value = 10.04058357522255
$ws.Range("D19").Value = 3.708355445821187
$ws.Range("F19").Value = 22.59415746919611
$ws.Range("G19").Value = 27.05476117932712
$ws.Range("H19").Value = 12.84780367604
$ws.Range("I19").Value = 17.88849294873842
$ws.Range("M19").Value = 21.8040042025756
$ws.Range("N19").Value = 16.22679928803851
$ws.Range("B20").Value = 12.61042699212446
$ws.Range("C20").Value = 10.19998096480379
$ws.Range("D20").Value = 3.714435447976215
$ws.Range("F20").Value = 22.68556879378328
$ws.Range("G20").Value = 27.24642555416386
$ws.Range("H20").Value = 12.85061249796365
$ws.Range("I20").Value = 17.88224188912119
$ws.Range("M20").Value = 22.04812161398913
$ws.Range("N20").Value = 16.20055435357183
$ws.Range("B21").Value = 13.25011595838464
$ws.Range("C21").Value = 10.71636794811675
$ws.Range("D21").Value = 3.734647075571485
$ws.Range("F21").Value = 22.99932280999699
$ws.Range("G21").Value = 27.89378798668896
$ws.Range("H21").Value = 12.86387034851905
$ws.Range("I21").Value = 17.86754691968971
$ws.Range("M21").Value = 22.8586213161173
$ws.Range("N21").Value = 16.11447366407124
$ws.Range("B22").Value = 13.69804232030823
$ws.Range("C22").Value = 11.03948915278182
$ws.Range("D22").Value = 3.74768811036299
$ws.Range("F22").Value = 23.20897840726101
$ws.Range("G22").Value = 28.31854400239757
$ws.Range("H22").Value = 12.87538032678898
$ws.Range("I22").Value = 17.86266669626738
$ws.Range("M22").Value = 23.38021060382656
$ws.Range("N22").Value = 16.05976840605321
$ws.Range("B23").Value = 13.45994928590924
$ws.Range("C23").Value = 10.86835247513351
$ws.Range("D23").Value = 3.740743735956852
$ws.Range("F23").Value = 23.09668236311117
$ws.Range("G23").Value = 28.09174573125525
$ws.Range("H23").Value = 12.86897712293451
$ws.Range("I23").Value = 17.86483778257431
$ws.Range("M23").Value = 23.10262561849114
$ws.Range("N23").Value = 16.08882418200498
$ws.Range("B24").Value = 12.59971151205773
$ws.Range("C24").Value = 10.19129948373426
$ws.Range("D24").Value = 3.714102358835749
$ws.Range("F24").Value = 22.68052314768864
$ws.Range("G24").Value = 27.23588524854782
$ws.Range("H24").Value = 12.85044387483117
$ws.Range("I24").Value = 17.88256149884273
$ws.Range("M24").Value = 22.03474996411292
$ws.Range("N24").Value = 16.2019876773335
$ws.Range("B25").Value = 11.62919647341128
$ws.Range("C25").Value = 9.399902689257967
$ws.Range("D25").Value = 3.684688783057232
$ws.Range("F25").Value = 22.2535418707564
$ws.Range("G25").Value = 26.32551707877973
$ws.Range("H25").Value = 12.8426773123957
$ws.Range("I25").Value = 17.92180677657364
$ws.Range("M25").Value = 20.85400412140049
$ws.Range("N25").Value = 16.33079223118141
